$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 2709
$ws.Range("F5").Value = 935
$ws.Range("F7").Value = 2275
$ws.Range("F8").Value = 1833
$ws.Range("F9").Value = 218
$ws.Range("F11").Value = 2479
$ws.Range("F12").Value = 550
$ws.Range("F13").Value = 242
$ws.Range("F14").Value = 55
$ws.Range("F16").Value = 128
$ws.Range("F17").Value = 115
$ws.Range("F18").Value = 9219
$ws.Range("F19").Value = 57
$ws.Range("F20").Value = 7158
$ws.Range("F21").Value = 11712
$ws.Range("F24").Value = 234
$ws.Range("F25").Value = 352
$ws.Range("F26").Value = 560
$ws.Range("F27").Value = 2589
$ws.Range("F28").Value = 233
$ws.Range("F29").Value = 197
$ws.Range("F30").Value = 2536
$ws.Range("F31").Value = 711
$ws.Range("F32").Value = 47
$ws.Range("F33").Value = 4510
$ws.Range("F34").Value = 906
$ws.Range("F35").Value = 354
$ws.Range("F37").Value = 527

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 1184
$ws.Range("F15").Value = 14
$ws.Range("F23").Value = 3

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 155

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 2709
$ws.Range("F7").Value = 935
$ws.Range("F9").Value = 2275
$ws.Range("F11").Value = 1834
$ws.Range("F13").Value = 218
$ws.Range("F14").Value = 2479
$ws.Range("F16").Value = 550
$ws.Range("F17").Value = 242
$ws.Range("F18").Value = 55
$ws.Range("F20").Value = 128
$ws.Range("F21").Value = 115
$ws.Range("F22").Value = 9219
$ws.Range("F23").Value = 57
$ws.Range("F24").Value = 7158
$ws.Range("F25").Value = 11712
$ws.Range("F28").Value = 234
$ws.Range("F29").Value = 352
$ws.Range("F31").Value = 560
$ws.Range("F33").Value = 2589
$ws.Range("F35").Value = 14
$ws.Range("F36").Value = 233
$ws.Range("F37").Value = 197
$ws.Range("F38").Value = 47
$ws.Range("F39").Value = 4510
$ws.Range("F46").Value = 527
$ws.Range("F48").Value = 3
